# ---------------------------------------------------------------------------
# Applies the commit "add delete/update in multi-partition table and release
# negative cases about array" to the mysql_updatedelete_cases workbook.
#
# Part 1: rows 66-82 (Testable column B) flip from "n" to "y" -- these are
#          the previously-disabled "array field" negative test cases that
#          are now being released.
# Part 2: four brand new test rows (100-103) are appended describing
#          multi-partition delete/update cases against schema9.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Part 1 - release the array negative cases: set Testable (column B) to "y"
# for rows 66 through 82 inclusive.
# ---------------------------------------------------------------------------
for ($r = 66; $r -le 82; $r++) {
    $ws.Range("B$r").Value = "y"
}

# ---------------------------------------------------------------------------
# Part 2 - append the four new multi-partition delete/update rows.
# Each new row is produced by duplicating the last existing data row (99),
# which carries the correct "text" number format / styling for this table,
# and then overwriting the cell values with the new test data.
# ---------------------------------------------------------------------------

# NOTE: this PowerShell engine does not reliably bind *named* parameters
# (e.g. "-RowNum 5"), so Add-TestRow takes its arguments positionally.
function Add-TestRow {
    param(
        [int]$RowNum,
        [string]$TestId,
        [string]$Title,
        [string]$Component,
        [string]$SchemaRef,
        [string]$ValueRef,
        [string]$Sql,
        [string]$EffectedRows,
        [string]$ValidationSql,
        [string]$ExpectedCsv,
        [string]$ValidationType
    )

    $ws.Rows("99").Copy()
    $ws.Rows("$RowNum").Insert(-4121)   # xlShiftDown

    $ws.Range("A$RowNum").Value = $TestId
    $ws.Range("B$RowNum").Value = "y"
    $ws.Range("C$RowNum").Value = $Title
    $ws.Range("D$RowNum").Value = $Component
    $ws.Range("E$RowNum").NumberFormat = "@"
    $ws.Range("F$RowNum").Value = $SchemaRef
    $ws.Range("G$RowNum").Value = $ValueRef
    $ws.Range("H$RowNum").Value = $Sql
    $ws.Range("I$RowNum").Value = $EffectedRows
    $ws.Range("J$RowNum").Value = $ValidationSql
    $ws.Range("K$RowNum").Value = $ExpectedCsv
    $ws.Range("L$RowNum").Value = $ValidationType
}

Add-TestRow 100 `
    "updel_099" `
    "多分区条件删除，小于条件" `
    "SQLFunction" `
    "schema9" `
    "updel_value04" `
    "delete from `$schema9 where id<100" `
    "9" `
    "select * from `$schema9" `
    "src/test/resources/io.dingodb.test/testdata/mysqlcases/dml/updatedelete/expectedresult/updatedelete_099.csv" `
    "csv_containsAll"

Add-TestRow 101 `
    "updel_100" `
    "多分区条件删除，大于条件" `
    "SQLFunction" `
    "schema9" `
    "updel_value04" `
    "delete from `$schema9 where id>0" `
    "10" `
    "select * from `$schema9" `
    "src/test/resources/io.dingodb.test/testdata/mysqlcases/dml/updatedelete/expectedresult/updatedelete_100.csv" `
    "csv_containsAll"

Add-TestRow 102 `
    "updel_101" `
    "多分区条件更新，小于条件" `
    "SQLFunction" `
    "schema9" `
    "updel_value04" `
    "update `$schema9 set name='Java' where id<20" `
    "6" `
    "select * from `$schema9" `
    "src/test/resources/io.dingodb.test/testdata/mysqlcases/dml/updatedelete/expectedresult/updatedelete_101.csv" `
    "csv_containsAll"

Add-TestRow 103 `
    "updel_102" `
    "多分区条件更新，大于条件" `
    "SQLFunction" `
    "schema9" `
    "updel_value04" `
    "update `$schema9 set amount=99.99 where id>=10" `
    "8" `
    "select * from `$schema9" `
    "src/test/resources/io.dingodb.test/testdata/mysqlcases/dml/updatedelete/expectedresult/updatedelete_102.csv" `
    "csv_containsAll"

# ---------------------------------------------------------------------------
# Update the sheet's recorded selection to match the post-edit state.
# ---------------------------------------------------------------------------
$ws.Range("B67").Select()

Write-Output "edit complete"
